# Friday Stop Loss Raise #1
#
# Updates a handful of "Stop Loss" related cells on Sheet1, moves the
# active selection, and switches the sheet's print orientation to
# portrait.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Raise stop-loss values (column M) for a few tickers, and bump the
# unrealized gain/loss figure in G7 to match the new stop loss.
$ws.Range("M3").Value = 97.11
$ws.Range("M5").Value = 66.79
$ws.Range("M6").Value = 45.55
$ws.Range("G7").Value = 67.07

# Switch the page setup to portrait orientation.
$ws.PageSetup.Orientation = 1

# Move/leave the active selection on O8, matching where the author's
# cursor ended up after making the edits.
$ws.Range("O8").Select()
